$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Delete the blank spacer row (old row 13) - shifts everything below up by one
#    and automatically shrinks the Total formula from SUM(J3:J13) to SUM(J3:J12).
$ws.Rows(13).Delete()

Write-Output "done"
